# Add season-record columns (Wins, Losses, Ties) to the worksheet.
# Mirrors the author's fix: the team's win/loss/tie record is appended
# as three new trailing columns (AD, AE, AF) for every row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should carry the same style as the existing header
# row (bold font + border + centered/top alignment). Copy that
# formatting from the last existing header cell (AC1) before writing
# the new header text, so the new cells reuse the same style record.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# 2008 Milwaukee Brewers season record: 90 wins, 72 losses, 0 ties.
$wins = 90
$losses = 72
$ties = 0

$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
